$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N ("Late"), shifting
# "Late" / heading("Outstanding" header) / "Outstanding" data right
# by one column. Copy column M's width onto the freshly inserted
# column so it keeps a fixed (non bestFit) width of 11 characters.
$mWidth = $ws.Columns("M").ColumnWidth
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $mWidth

# Make "Repayment schedule" the active sheet/tab, and leave the
# selection on K14 as last interacted cell.
$ws.Activate()
$ws.Range("K14").Select() | Out-Null
